$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("航天发展", "利欧股份", "利欧股份"),
    @("中国卫星", "航天电子", "三花智控"),
    @("利欧股份", "金风科技", "航天发展"),
    @("蓝色光标", "中国卫星", "金风科技"),
    @("雷科防务", "蓝色光标", "万向钱潮"),
    @("三花智控", "雷科防务", "平潭发展"),
    @("金风科技", "三花智控", "雷科防务"),
    @("神剑股份", "航天发展", "神剑股份"),
    @("顺灏股份", "中国卫通", "中国卫星"),
    @("航天电子", "神剑股份", "山子高科"),
    @("五洲新春", "洲际油气", "蓝色光标"),
    @("山子高科", "五洲新春", "五洲新春"),
    @("泰尔股份", "泰尔股份", "泰尔股份"),
    @("中国卫通", "山子高科", "航天电子"),
    @("万向钱潮", "北斗星通", "顺灏股份"),
    @("贵州茅台", "顺灏股份", "首开股份"),
    @("信维通信", "岩山科技", "卧龙电驱"),
    @("北斗星通", "烽火通信", "中国卫通"),
    @("福莱新材", "中国海油", "北斗星通"),
    @("洲际油气", "王子新材", "大业股份")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $ws.Cells.Item($rowNum, 1).Value = $data[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $data[$i][1]
    $ws.Cells.Item($rowNum, 3).Value = $data[$i][2]
}

$wb.Save()
